$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New FedEx tracking numbers (ShipmentTracking, column P) for rows 2..26,
# replacing the previous batch of tracking numbers one-for-one.
$newTrackingNumbers = @(
    "320017962708",
    "320017962719",
    "320017962741",
    "320017962774",
    "320017962811",
    "320017962833",
    "320017962866",
    "320017962888",
    "320017962936",
    "320017962958",
    "320017962991",
    "320017963016",
    "320017963049",
    "320017963060",
    "320017963093",
    "320017963119",
    "320017963152",
    "320017963174",
    "320017963200",
    "320017963222",
    "320017963255",
    "320017963266",
    "320017963288",
    "320017963299",
    "320017963314"
)

for ($i = 0; $i -lt $newTrackingNumbers.Length; $i++) {
    $row = $i + 2
    $cell = $ws.Cells.Item($row, 16)
    # Force text storage so the all-digit tracking number isn't coerced to a
    # number, then restore the default "Normal" style so no stray cell
    # format/style is left behind.
    $cell.NumberFormat = "@"
    $cell.Value2 = $newTrackingNumbers[$i]
    $cell.Style = "Normal"
}
